# Generate Report for Handoff
# Updates the localization-status report with a fresh handoff run:
#   - new GUID-based file id (89648d50-... -> 3ed6cbfc-...)
#   - new hash for the generated .xlf targets (8fa999f7... -> 84f60d57...)
#   - refreshed handoff timestamps

$wb = $excel.ActiveWorkbook

$newId   = "3ed6cbfc-9d66-4dc2-b7ef-0658311d60c8"
$newHash = "84f60d574527c0f3aa39a14fa60b3a6552f0e879"

$newMd      = "$newId.md"
$newZhXlf   = "$newId.$newHash.zh-cn.xlf"
$newDeXlf   = "$newId.$newHash.de-de.xlf"

$newHandoffDate     = "2016-03-21 17:01:15"
$newZhHandoffStamp  = "2016-03-21 17:01:08"

function Set-HyperlinkDisplay($ws, $cellRef, $newText) {
    $target = $ws.Range($cellRef).Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            $hl.TextToDisplay = $newText
        }
    }
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = $newHandoffDate
Set-HyperlinkDisplay $wsOverview "A2" $newMd

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhHandoffStamp
Set-HyperlinkDisplay $wsZh "A2" $newMd
Set-HyperlinkDisplay $wsZh "D2" $newZhXlf

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newHandoffDate
Set-HyperlinkDisplay $wsDe "A2" $newMd
Set-HyperlinkDisplay $wsDe "D2" $newDeXlf
